$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "aliases" column (C) used to be computed via a formula (=D2 etc.) that
# referenced the "aliasesx" column (D), which itself held the (buggy) result of
# a regex used to build the alias list. The fix: replace the formula results in
# column C with the literal values from column D (the corrected values), then
# delete column D entirely (shifting sequence/accession left).

$lastRow = $ws.Cells(1, 1).End(-4121).Row   # xlDown
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 4).Value2
}

# Delete entire column D (old "aliasesx"). This shifts E->D, F->E.
$ws.Columns.Item(4).Delete()

# Best-fit the gene-name column (B) and the newly-literal aliases column (C)
# to their content, matching the widths Excel computed after the edit.
$ws.Columns.Item(2).ColumnWidth = 9.1640625
$ws.Columns.Item(3).ColumnWidth = 182

# Update selection to match the post-edit cursor position recorded in the diff.
$ws.Range("A3").Select()
